# Refresh the cryptos list numbers (price + 1h volume/change columns).
#
# Column D ("Price") and column E ("Volume(1h)") are stored as plain text
# in this sheet (prices use a dotted/European-style grouping like
# "41.524.34", and percentages carry padding spaces like "  +1.09%  "),
# so every one of these must stay a literal string, not get silently
# re-interpreted as a number/percentage by Excel.
#
# Most of the new price strings (multi-dot, e.g. "41.525.48") can never
# parse as a number, so a normal .Value assignment already keeps them as
# text. A few new prices (e.g. "313.71") DO look like plain numbers, so
# they are written with a leading apostrophe - exactly like a user typing
# '313.71 into a cell - which forces Excel to keep them as text instead of
# converting them to a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    # Leading apostrophe = Excel "treat as text" quote-prefix, same as
    # typing it interactively. Excel strips the apostrophe itself and
    # keeps it out of the stored string.
    $ws.Range($Cell).Value = "'" + $Text
}

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "41.525.48"
$ws.Range("E2").Value = "  +1.40%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "2.478.59"
$ws.Range("E3").Value = "  +1.29%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  -0.10%  "

# --- Row 5 (BNB) ---
Set-TextValue "D5" "313.71"
$ws.Range("E5").Value = "  +1.49%  "

# --- Row 6 (Solana) ---
Set-TextValue "D6" "93.25"
$ws.Range("E6").Value = "  +1.64%  "

# --- Row 7 (XRP) ---
$ws.Range("E7").Value = "  -1.11%  "

# --- Row 8 (USDC) ---
$ws.Range("E8").Value = "  -0.15%  "

# --- Row 9 (Cardano) ---
$ws.Range("E9").Value = "  +3.24%  "

# --- Row 10 (Avalanche) ---
Set-TextValue "D10" "32.72"
$ws.Range("E10").Value = "  -0.07%  "

# --- Row 11 (Dogecoin) ---
$ws.Range("E11").Value = "  +2.13%  "

# --- Row 12 (TRON) ---
$ws.Range("E12").Value = "  +3.15%  "

# --- Row 13 (Wrapped liquid staked Ether 2.0) ---
$ws.Range("D13").Value = "2.864.64"
$ws.Range("E13").Value = "  +1.25%  "

# --- Row 14 (Polkadot) ---
Set-TextValue "D14" "6.83"
$ws.Range("E14").Value = "  -0.45%  "

# --- Row 15 (Chainlink) ---
Set-TextValue "D15" "16.11"
$ws.Range("E15").Value = "  +10.47%  "

# --- Row 16 (Wrapped Ether) ---
$ws.Range("D16").Value = "2.515.60"
$ws.Range("E16").Value = "  +2.89%  "

# --- Row 17 (Polygon) ---
$ws.Range("E17").Value = "  -0.74%  "

# --- Row 18 (Wrapped BTC) ---
$ws.Range("D18").Value = "41.518.31"
$ws.Range("E18").Value = "  +1.41%  "

# --- Row 19 (Uniswap) ---
$ws.Range("E19").Value = "  +3.44%  "

# --- Row 20 (Shiba Inu) ---
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +3.77%  "

# --- Row 21 (Litecoin) ---
Set-TextValue "D21" "71.74"
$ws.Range("E21").Value = "  +6.52%  "

# --- Row 22 (Internet Computer / DFINITY) ---
Set-TextValue "D22" "11.35"
$ws.Range("E22").Value = "  +3.98%  "

# --- Row 23 (Bitcoin Cash) ---
Set-TextValue "D23" "236.49"
$ws.Range("E23").Value = "  +1.66%  "

# --- Row 24 (PancakeSwap) ---
Set-TextValue "D24" "2.71"
$ws.Range("E24").Value = "  -0.60%  "

# --- Row 25 (Dai) ---
$ws.Range("E25").Value = "  -0.27%  "

# --- Row 26 (ImmutableX) ---
$ws.Range("E26").Value = "  +1.56%  "

# --- Row 27 (Ethereum Classic) ---
Set-TextValue "D27" "24.87"
$ws.Range("E27").Value = "  +5.97%  "

# --- Row 28 (Toncoin) ---
$ws.Range("E28").Value = "  +0.63%  "

# --- Row 29 (Cosmos) ---
Set-TextValue "D29" "9.64"
$ws.Range("E29").Value = "  +2.09%  "

# --- Row 30 (Injective Protocol) ---
Set-TextValue "D30" "35.86"
$ws.Range("E30").Value = "  +2.05%  "

# --- Row 31 (Monero) ---
Set-TextValue "D31" "157.99"
$ws.Range("E31").Value = "  +5.24%  "

# --- Row 32 (Filecoin) ---
$ws.Range("E32").Value = "  +1.98%  "

# --- Row 33 (WEMIX Token) ---
$ws.Range("E33").Value = "  +1.56%  "

# --- Row 34 (Hedera) ---
$ws.Range("E34").Value = "  +3.76%  "

# --- Row 35 (Celestia) ---
Set-TextValue "D35" "17.40"
$ws.Range("E35").Value = "  +5.68%  "

# --- Row 36 (ApeX Protocol) ---
$ws.Range("E36").Value = "  -7.85%  "

# --- Row 37 (Lido DAO Token) ---
$ws.Range("E37").Value = "  -0.72%  "

# --- Row 38 (Kaspa) ---
$ws.Range("E38").Value = "  +4.47%  "

# --- Row 39 (ARBITRUM) ---
$ws.Range("E39").Value = "  +0.22%  "

# --- Row 40 (Stellar) ---
$ws.Range("E40").Value = "  +1.05%  "

# --- Row 41 (Render Token) ---
$ws.Range("E41").Value = "  +0.67%  "

# --- Row 42 (First Digital USD) ---
$ws.Range("E42").Value = "  -0.11%  "

# --- Row 43 (EnergySwap) ---
Set-TextValue "D43" "19.71"
$ws.Range("E43").Value = "  +1.19%  "

# --- Row 44 (Maker) ---
$ws.Range("D44").Value = "1.970.43"
$ws.Range("E44").Value = "  +1.00%  "

# --- Row 45 (VeChain) ---
$ws.Range("E45").Value = "  +1.79%  "

# --- Row 46 (NEAR Protocol) ---
Set-TextValue "D46" "2.95"
$ws.Range("E46").Value = "  -0.31%  "

# --- Row 47 (Frax Share) ---
Set-TextValue "D47" "9.15"
$ws.Range("E47").Value = "  +8.14%  "

# --- Row 48 (Rocket Pool ETH) ---
$ws.Range("D48").Value = "2.722.44"
$ws.Range("E48").Value = "  +1.33%  "

# --- Row 49 (Aave) ---
Set-TextValue "D49" "97.78"
$ws.Range("E49").Value = "  +2.78%  "

# --- Row 50 (ordi) ---
Set-TextValue "D50" "68.03"
$ws.Range("E50").Value = "  -1.27%  "

# --- Row 51 (BitcoinSV) ---
Set-TextValue "D51" "72.34"
$ws.Range("E51").Value = "  -0.55%  "
